$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.228395223617554
$ws.Range("B1").Value = 2.722322463989258
$ws.Range("C1").Value = 2.87277364730835
$ws.Range("D1").Value = 2.543979406356812
$ws.Range("E1").Value = 0.8329096436500549
